# Fruta / hortaliza, semanal
# The weekly refresh reorders the historical price rows (2-10): each row's
# data set (date, variety, quality, volume, prices, unit, origin, etc.)
# moves to a new row position. Columns A,B,C,E,F,G,H,I,J stay constant for
# every data row, so only D and K:T need to be relocated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that are permuted per-row.
$cols = @("D","K","L","M","N","O","P","Q","R","S","T")

# Snapshot the current ("before") values for rows 2..10 across the columns
# that vary, before any writes happen (so sources aren't clobbered mid-way).
$snapshot = @{}
for ($r = 2; $r -le 10; $r++) {
    $rowVals = @{}
    foreach ($col in $cols) {
        $rowVals[$col] = $ws.Range("$col$r").Value()
    }
    $snapshot[$r] = $rowVals
}

# New row r receives the old data that used to live in row $map[r].
$map = @{
    2  = 4
    3  = 5
    4  = 6
    5  = 3
    6  = 7
    7  = 8
    8  = 9
    9  = 10
    10 = 2
}

foreach ($r in $map.Keys) {
    $src = $map[$r]
    $rowVals = $snapshot[$src]
    foreach ($col in $cols) {
        $ws.Range("$col$r").Value = $rowVals[$col]
    }
}
